$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the Price cells that are about to receive a new numeric-looking
# string (e.g. "214.70", "0.996") as Text first, so Excel keeps them as
# literal strings instead of silently converting them to floats/doubles.
$ws.Range("D2,D3,D4,D5,D6,D10,D11,D12,D13,D15,D16,D17,D18,D20,D21,D23,D24,D25,D26,D27,D29,D30,D31,D33,D36,D37,D38,D41,D42,D43,D44,D46,D47,D49,D50,D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.707.46"
$ws.Range("E2").Value = "  +1.72%  "

$ws.Range("D3").Value = "1.624.74"
$ws.Range("E3").Value = "  +2.17%  "

$ws.Range("D4").Value = "0.996"
$ws.Range("E4").Value = "  -0.43%  "

$ws.Range("D5").Value = "214.70"
$ws.Range("E5").Value = "  +1.14%  "

$ws.Range("D6").Value = "0.505"
$ws.Range("E6").Value = "  +0.84%  "

$ws.Range("E7").Value = "  -0.34%  "

$ws.Range("E8").Value = "  +0.65%  "

$ws.Range("E9").Value = "  +0.58%  "

$ws.Range("D10").Value = "19.37"
$ws.Range("E10").Value = "  -0.07%  "

$ws.Range("D11").Value = "0.0855"
$ws.Range("E11").Value = "  +0.73%  "

$ws.Range("D12").Value = "1.851.86"
$ws.Range("E12").Value = "  +2.12%  "

$ws.Range("D13").Value = "1.624.84"
$ws.Range("E13").Value = "  +1.92%  "

$ws.Range("E14").Value = "  +0.37%  "

$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "65.13"
$ws.Range("E15").Value = "  +1.18%  "

$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "0.514"
$ws.Range("E16").Value = "  -1.28%  "

$ws.Range("D17").Value = "26.705.07"
$ws.Range("E17").Value = "  +1.72%  "

$ws.Range("D18").Value = "231.68"
$ws.Range("E18").Value = "  +8.43%  "

$ws.Range("E19").Value = "  +0.22%  "

$ws.Range("D20").Value = "7.64"
$ws.Range("E20").Value = "  +2.95%  "

$ws.Range("D21").Value = "0.997"
$ws.Range("E21").Value = "  -0.34%  "

$ws.Range("E22").Value = "  +2.45%  "

$ws.Range("D23").Value = "2.23"
$ws.Range("E23").Value = "  +4.01%  "

$ws.Range("D24").Value = "9.12"
$ws.Range("E24").Value = "  +1.12%  "

$ws.Range("D25").Value = "145.46"
$ws.Range("E25").Value = "  +1.23%  "

$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.36%  "

$ws.Range("D27").Value = "7.04"
$ws.Range("E27").Value = "  -0.21%  "

$ws.Range("E28").Value = "  +2.36%  "

$ws.Range("D29").Value = "15.66"
$ws.Range("E29").Value = "  +2.96%  "

$ws.Range("D30").Value = "0.0499"
$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("D31").Value = "1.17"
$ws.Range("E31").Value = "  +0.72%  "

$ws.Range("E32").Value = "  +1.79%  "

$ws.Range("D33").Value = "1.447.78"
$ws.Range("E33").Value = "  +8.03%  "

$ws.Range("E34").Value = "  +2.33%  "

$ws.Range("D36").Value = "1.49"
$ws.Range("E36").Value = "  +0.99%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.0168"
$ws.Range("E37").Value = "  +0.81%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "0.560"
$ws.Range("E38").Value = "  -5.27%  "

$ws.Range("E39").Value = "  +2.80%  "

$ws.Range("E40").Value = "  +1.96%  "

$ws.Range("D41").Value = "0.997"
$ws.Range("E41").Value = "  -0.37%  "

$ws.Range("D42").Value = "2.29"
$ws.Range("E42").Value = "  +6.68%  "

$ws.Range("D43").Value = "0.951"
$ws.Range("E43").Value = "  -5.26%  "

$ws.Range("D44").Value = "1.763.73"
$ws.Range("E44").Value = "  +2.25%  "

$ws.Range("E45").Value = "  -0.46%  "

$ws.Range("D46").Value = "62.20"
$ws.Range("E46").Value = "  +0.41%  "

$ws.Range("D47").Value = "88.62"
$ws.Range("E47").Value = "  +3.31%  "

$ws.Range("E48").Value = "  +2.47%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0504"
$ws.Range("E49").Value = "  +0.34%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.0968"
$ws.Range("E50").Value = "  -0.86%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.52"
$ws.Range("E51").Value = "  +1.86%  "
